# Add new columns I (I0) and J (IF) to Sheet1, mirroring the styling of the
# existing header/data columns, and populate the values for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style used by the other header cells (e.g. H1) onto the new
# header cells so formatting (bold, borders, centered) matches.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data rows (rows 2-36) for column I (I0) and column J (IF) ---
$values = @{
    2  = @(8, 8)
    3  = @(6, 6)
    4  = @(2, 3)
    5  = @(9, 9)
    6  = @(6, 7)
    7  = @(6, 7)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(7, 8)
    11 = @(8, 8)
    12 = @(7, 7)
    13 = @(6, 6)
    14 = @(5, 5)
    15 = @(6, 7)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(7, 8)
    19 = @(12, 12)
    20 = @(8, 9)
    21 = @(8, 8)
    22 = @(7, 8)
    23 = @(7, 7)
    24 = @(6, 6)
    25 = @(8, 8)
    26 = @(7, 7)
    27 = @(8, 8)
    28 = @(5, 6)
    29 = @(5, 5)
    30 = @(5, 5)
    31 = @(1, 1)
    32 = @(6, 6)
    33 = @(6, 6)
    34 = @(6, 6)
    35 = @(5, 5)
    36 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
